# Apply scheduled profit-sheet recalculation updates across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 501.1
$ws.Range("I38").Value = 14.571428
$ws.Range("K38").Value = 43.714284
$ws.Range("M38").Value = 328.285716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2121.1516
$ws.Range("I61").Value = 1597.8889
$ws.Range("J61").Value = 4475.8335
$ws.Range("K61").Value = 1597.8889
$ws.Range("L61").Value = 4475.8335
$ws.Range("M61").Value = -1385.8889
$ws.Range("N61").Value = -4899.8335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3481.647
$ws.Range("I74").Value = 970.2941
$ws.Range("J74").Value = 8504.352999999999
$ws.Range("K74").Value = 970.2941
$ws.Range("L74").Value = 8504.352999999999
$ws.Range("M74").Value = -96.29409999999996
$ws.Range("N74").Value = -10252.353

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3481.647
$ws.Range("I77").Value = 970.2941
$ws.Range("J77").Value = 8504.352999999999
$ws.Range("K77").Value = 4851.470499999999
$ws.Range("L77").Value = 42521.765
$ws.Range("M77").Value = -483.4704999999994
$ws.Range("N77").Value = -51257.765

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2137.9375
$ws.Range("I132").Value = 1695.3137
$ws.Range("J132").Value = 3874.3845
$ws.Range("K132").Value = 5085.9411
$ws.Range("L132").Value = 11623.1535
$ws.Range("M132").Value = -2555.9411
$ws.Range("N132").Value = -16683.1535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2121.1516
$ws.Range("I136").Value = 1597.8889
$ws.Range("J136").Value = 4475.8335
$ws.Range("K136").Value = 4793.6667
$ws.Range("L136").Value = 13427.5005
$ws.Range("M136").Value = -2243.6667
$ws.Range("N136").Value = -18527.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1591.4615
$ws.Range("I107").Value = 1470.6364
$ws.Range("J107").Value = 2256
$ws.Range("K107").Value = 1470.6364
$ws.Range("L107").Value = 2256
$ws.Range("M107").Value = 449.3635999999999
$ws.Range("N107").Value = -6096

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 41612.31
$ws.Range("J132").Value = 41612.31
$ws.Range("L132").Value = 41612.31
$ws.Range("N132").Value = -51732.31

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 76925700
$ws.Range("J4").Value = 83336080
$ws.Range("L4").Value = 83336080
$ws.Range("N4").Value = -83336304

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4370.0625
$ws.Range("I31").Value = 2150.889
$ws.Range("J31").Value = 7223.2856
$ws.Range("K31").Value = 2150.889
$ws.Range("L31").Value = 7223.2856
$ws.Range("M31").Value = -1855.889
$ws.Range("N31").Value = -7813.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4370.0625
$ws.Range("I34").Value = 2150.889
$ws.Range("J34").Value = 7223.2856
$ws.Range("K34").Value = 2150.889
$ws.Range("L34").Value = 7223.2856
$ws.Range("M34").Value = -1948.889
$ws.Range("N34").Value = -7627.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1665.5312
$ws.Range("I58").Value = 1281.3158
$ws.Range("J58").Value = 2227.077
$ws.Range("K58").Value = 1281.3158
$ws.Range("L58").Value = 2227.077
$ws.Range("M58").Value = -1078.3158
$ws.Range("N58").Value = -2633.077

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1520.6511
$ws.Range("I132").Value = 1038.8379
$ws.Range("K132").Value = 3116.5137
$ws.Range("M132").Value = -586.5137

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2589.5642
$ws.Range("I134").Value = 1481.1666
$ws.Range("J134").Value = 6284.222
$ws.Range("K134").Value = 4443.4998
$ws.Range("L134").Value = 18852.666
$ws.Range("M134").Value = -1908.4998
$ws.Range("N134").Value = -23922.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1665.5312
$ws.Range("I136").Value = 1281.3158
$ws.Range("J136").Value = 2227.077
$ws.Range("K136").Value = 3843.9474
$ws.Range("L136").Value = 6681.231000000001
$ws.Range("M136").Value = -1293.9474
$ws.Range("N136").Value = -11781.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4976499
$ws.Range("I131").Value = 408.66666
$ws.Range("J131").Value = 5748651
$ws.Range("K131").Value = 1225.99998
$ws.Range("L131").Value = 17245953
$ws.Range("M131").Value = 3814.00002
$ws.Range("N131").Value = -17256033

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 6980
$ws.Range("J44").Value = 6980
$ws.Range("L44").Value = 6980
$ws.Range("N44").Value = -8172

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8990
$ws.Range("I70").Value = 10128.571
$ws.Range("J70").Value = 6333.3335
$ws.Range("K70").Value = 10128.571
$ws.Range("L70").Value = 6333.3335
$ws.Range("M70").Value = -9858.571
$ws.Range("N70").Value = -6873.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8990
$ws.Range("I73").Value = 10128.571
$ws.Range("J73").Value = 6333.3335
$ws.Range("K73").Value = 10128.571
$ws.Range("L73").Value = 6333.3335
$ws.Range("M73").Value = -9192.571
$ws.Range("N73").Value = -8205.333500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2835.7144
$ws.Range("I80").Value = 2700
$ws.Range("J80").Value = 3333.3333
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 3333.3333
$ws.Range("M80").Value = -1702
$ws.Range("N80").Value = -5329.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2835.7144
$ws.Range("I83").Value = 2700
$ws.Range("J83").Value = 3333.3333
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 16666.6665
$ws.Range("M83").Value = -8508
$ws.Range("N83").Value = -26650.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 11205.6
$ws.Range("J123").Value = 11205.6
$ws.Range("L123").Value = 11205.6
$ws.Range("N123").Value = -16105.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 19354.625
$ws.Range("J136").Value = 19163.467
$ws.Range("L136").Value = 57490.401
$ws.Range("N136").Value = -62590.401

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3237
$ws.Range("I7").Value = 2125.75
$ws.Range("K7").Value = 2125.75
$ws.Range("M7").Value = -2013.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2396.1428
$ws.Range("I68").Value = 1992.5
$ws.Range("J68").Value = 2934.3333
$ws.Range("K68").Value = 1992.5
$ws.Range("L68").Value = 2934.3333
$ws.Range("M68").Value = -1243.5
$ws.Range("N68").Value = -4432.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2396.1428
$ws.Range("I71").Value = 1992.5
$ws.Range("J71").Value = 2934.3333
$ws.Range("K71").Value = 9962.5
$ws.Range("L71").Value = 14671.6665
$ws.Range("M71").Value = -6218.5
$ws.Range("N71").Value = -22159.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3908
$ws.Range("I122").Value = 2875
$ws.Range("J122").Value = 4022.7778
$ws.Range("K122").Value = 8625
$ws.Range("L122").Value = 12068.3334
$ws.Range("M122").Value = -6175
$ws.Range("N122").Value = -16968.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3237
$ws.Range("I126").Value = 2125.75
$ws.Range("K126").Value = 6377.25
$ws.Range("M126").Value = -3907.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 648.93335
$ws.Range("I113").Value = 648.93335
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1946.80005
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 223.1999499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1363.9445
$ws.Range("I122").Value = 1153.125
$ws.Range("J122").Value = 1785.5834
$ws.Range("K122").Value = 3459.375
$ws.Range("L122").Value = 5356.7502
$ws.Range("M122").Value = -1009.375
$ws.Range("N122").Value = -10256.7502

# Row 113 on WVR no longer carries an HQ-profit figure; clear N113 outright.
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N113").ClearContents()
